# Apply row-data permutation to the "red_flag_10" procuring-entity table.
# The underlying records for rows 2, 3, 4, 5, 6 and 8 are swapped among
# themselves (row 7 - INTERSEC SECURITY COMPANY - is untouched). Only the
# cells whose value actually differs between the old and new row content
# are written, and a cell that has no value under its new row is cleared.
#
# Columns K (contact_point_telephone), N (contact_person_nat_id) and
# Q (ext_tin) hold numeric-looking text (phone numbers with leading
# zeros, long national-id/TIN numbers) that Excel would otherwise
# auto-coerce into numbers, losing leading zeros / precision, so those
# cells are explicitly formatted as text before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '3114072400000X'
$ws.Range("B2").Value = 'CONNECT GLOBAL BUSINESS COMPANY Ltd'
$ws.Range("D2").Value = 'CONNECT GLOBAL BUSINESS COMPANY Ltd'
$ws.Range("E2").Value = 'KIGALI - NYARUGENGE'
$ws.Range("F2").Value = 'KIGALI - NYARUGENGE'
$ws.Range("G2").Value = 'KIGALI - NYARUGENGE'
$ws.Range("H2").Value = 'POB:5564 Kigali'
$ws.Range("J2").Value = 'niyafeos@yahoo.fr'
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = '788757320'
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = '1197980004514176'
$ws.Range("O2").Value = 'RWANDA'
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = '106907027'

# Row 3
$ws.Range("A3").Value = '3100580400000V'
$ws.Range("B3").Value = 'SOCIETE RWANDAISE DES PNEUMATIQUES BANDAG  LTD'
$ws.Range("D3").Value = 'SOCIETE RWANDAISE DES PNEUMATIQUES BANDAG  LTD'
$ws.Range("E3").Value = 'Kicukiro- Rwanda'
$ws.Range("F3").Value = 'Kicukiro- Rwanda'
$ws.Range("G3").Value = 'Kicukiro- Rwanda'
$ws.Range("H3").Value = '1132 kigali'
$ws.Range("J3").Value = 'bandag@subizo.com'
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = '0788303361'
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = 'AA3035975'
$ws.Range("O3").Value = 'ITALY'
$ws.Range("P3").Value = 'I&M Bank'
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = '100003368'

# Row 4
$ws.Range("A4").Value = '3109596400000L'
$ws.Range("B4").Value = 'NEPOMSCENE BUSINESS COMPANY  LTD'
$ws.Range("D4").Value = 'NEPOMSCENE BUSINESS COMPANY  LTD'
$ws.Range("E4").Value = 'NYAMAGABE District,Gasaka Sector'
$ws.Range("F4").Value = 'NYAMAGABE District,Gasaka Sector'
$ws.Range("G4").Value = 'NYAMAGABE District,Gasaka Sector'
$ws.Range("H4").Value = $null
$ws.Range("J4").Value = 'nbchvgmn@gmail.com'
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = '0788475217'
$ws.Range("M4").Value = $null
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = '1198380172943184'
$ws.Range("O4").Value = 'RWANDA'
$ws.Range("P4").Value = 'EQUITY BANK'
$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = '103496546'

# Row 5
$ws.Range("A5").Value = '3100989200000S'
$ws.Range("B5").Value = 'H.VEDASTE COMPANY Ltd'
$ws.Range("D5").Value = 'H.VEDASTE COMPANY Ltd'
$ws.Range("E5").Value = 'Kicukiro/Kigali'
$ws.Range("F5").Value = 'Kicukiro/Kigali'
$ws.Range("G5").Value = 'Kicukiro/Kigali'
$ws.Range("H5").Value = '25 NGOMA'
$ws.Range("J5").Value = 'hitimanaveda16@gmail.com'
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = '0788562686'
$ws.Range("M5").Value = $null
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = '1196280064731068'
$ws.Range("P5").Value = 'BPR'
$ws.Range("Q5").NumberFormat = "@"
$ws.Range("Q5").Value = '115405774'

# Row 6
$ws.Range("A6").Value = '3102226900000M'
$ws.Range("B6").Value = 'BROADBAND SYSTEMS CORPORATION  LTD'
$ws.Range("D6").Value = 'BROADBAND SYSTEMS CORPORATION  LTD'
$ws.Range("E6").Value = 'Remera, Gisimenti Airport Road (kn5 Rda), Opposite ChezLando'
$ws.Range("F6").Value = 'Remera, Gisimenti Airport Road (kn5 Rda), Opposite ChezLando'
$ws.Range("G6").Value = 'Remera, Gisimenti Airport Road (kn5 Rda), Opposite ChezLando'
$ws.Range("H6").Value = '7229 KIGALI, RWANDA'
$ws.Range("J6").Value = 'gilbert.kayinamura@bsc.rw'
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = '4141'
$ws.Range("M6").Value = 'www.bsc.rw'
$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value = '1198080006094070'
$ws.Range("P6").Value = 'BK'
$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = '101982714'

# Row 8
$ws.Range("A8").Value = '3100960800000N'
$ws.Range("B8").Value = 'INYANGE INDUSTRIES  LTD'
$ws.Range("D8").Value = 'INYANGE INDUSTRIES  LTD'
$ws.Range("E8").Value = 'MASAKA GASABO DISTRICT'
$ws.Range("F8").Value = 'MASAKA GASABO DISTRICT'
$ws.Range("G8").Value = 'MASAKA GASABO DISTRICT'
$ws.Range("H8").Value = '4584 kigali-rwanda'
$ws.Range("J8").Value = 'bjames@inyangeindustries.com'
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = '0788309662'
$ws.Range("M8").Value = 'www.inyangeindustries.com'
$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = '1197580006310064'
$ws.Range("O8").Value = 'Rwandan'
$ws.Range("P8").Value = 'NCBA'
$ws.Range("Q8").NumberFormat = "@"
$ws.Range("Q8").Value = '100095380'
